# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates1 = @{
    "F2"  = 250
    "F5"  = 434
    "F9"  = 6813
    "F16" = 16220
    "F17" = 1594
    "F19" = 331
    "F22" = 11378
    "F24" = 1019
    "F26" = 320
}

$updates4 = @{
    "F2"  = 250
    "F5"  = 434
    "F10" = 6813
    "F18" = 16220
    "F19" = 1594
    "F21" = 331
    "F26" = 11378
    "F28" = 1019
    "F30" = 320
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $updates1.Keys) {
    $ws1.Range($addr).Value = $updates1[$addr]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $updates4.Keys) {
    $ws4.Range($addr).Value = $updates4[$addr]
}
